# Updates the country/stat table on sheet "Pais" with the refreshed
# snapshot (new timestamp + refreshed case counts) and reorders a handful
# of countries whose totals changed rank: Peru/Colombia, Chile/Francia,
# China/Guatemala, Timor Oriental/Santa Lucia, Islas Malvinas/Montserrat.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 21 de Septiembre de 2020 a las 01:06'

$ws.Range("B4").Value = 7000193
$ws.Range("C4").Value = 31814
$ws.Range("D4").Value = 4245816
$ws.Range("E4").Value = 2550264
$ws.Range("G4").Value = 289
$ws.Range("H4").Value = 204113

$ws.Range("D6").Value = 3851227
$ws.Range("E6").Value = 556507

$ws.Range("A8").Value = 'Colombia'
$ws.Range("B8").Value = 765076
$ws.Range("C8").Value = 6678
$ws.Range("D8").Value = 633199
$ws.Range("E8").Value = 107669
$ws.Range("G8").Value = 169
$ws.Range("H8").Value = 24208

$ws.Range("A9").Value = 'Peru'
$ws.Range("B9").Value = 762865
$ws.Range("D9").Value = 607837
$ws.Range("E9").Value = 123659
$ws.Range("H9").Value = 31369

$ws.Range("B13").Value = 631365
$ws.Range("C13").Value = 8431
$ws.Range("E13").Value = 130081
$ws.Range("G13").Value = 254
$ws.Range("H13").Value = 13053

$ws.Range("A14").Value = 'Francia'
$ws.Range("B14").Value = 452763
$ws.Range("C14").Value = 10569
$ws.Range("D14").Value = 91574
$ws.Range("E14").Value = 329904
$ws.Range("G14").Value = 11
$ws.Range("H14").Value = 31285

$ws.Range("A15").Value = 'Chile'
$ws.Range("B15").Value = 446274
$ws.Range("C15").Value = 1600
$ws.Range("D15").Value = 419746
$ws.Range("E15").Value = 14242
$ws.Range("G15").Value = 32
$ws.Range("H15").Value = 12286

$ws.Range("B36").Value = 106203
$ws.Range("C36").Value = 602
$ws.Range("D36").Value = 81365
$ws.Range("E36").Value = 22581
$ws.Range("G36").Value = 10
$ws.Range("H36").Value = 2257

$ws.Range("B37").Value = 102015
$ws.Range("C37").Value = 115
$ws.Range("D37").Value = 89532
$ws.Range("E37").Value = 6713
$ws.Range("G37").Value = 20
$ws.Range("H37").Value = 5770

$ws.Range("A44").Value = 'Guatemala'
$ws.Range("B44").Value = 85444
$ws.Range("C44").Value = 292
$ws.Range("D44").Value = 74859
$ws.Range("E44").Value = 7466
$ws.Range("G44").Value = 14
$ws.Range("H44").Value = 3119

$ws.Range("A45").Value = 'China'
$ws.Range("B45").Value = 85279
$ws.Range("C45").Value = 10
$ws.Range("D45").Value = 80477
$ws.Range("E45").Value = 168
$ws.Range("H45").Value = 4634

$ws.Range("B48").Value = 78657
$ws.Range("C48").Value = 584
$ws.Range("D48").Value = 71030
$ws.Range("E48").Value = 6127
$ws.Range("G48").Value = 5
$ws.Range("H48").Value = 1500

$ws.Range("B53").Value = 66656
$ws.Range("C53").Value = 707
$ws.Range("D53").Value = 56096
$ws.Range("E53").Value = 10013
$ws.Range("G53").Value = 8
$ws.Range("H53").Value = 547

$ws.Range("B58").Value = 57242
$ws.Range("C58").Value = 97
$ws.Range("D58").Value = 48569
$ws.Range("E58").Value = 7575
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 1098

$ws.Range("B97").Value = 10325
$ws.Range("C97").Value = 39
$ws.Range("D97").Value = 9692
$ws.Range("E97").Value = 569
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 64

$ws.Range("B109").Value = 7368
$ws.Range("C109").Value = 3
$ws.Range("D109").Value = 6951
$ws.Range("E109").Value = 256

$ws.Range("B154").Value = 1917
$ws.Range("C154").Value = 13
$ws.Range("D154").Value = 1621
$ws.Range("E154").Value = 250

$ws.Range("B165").Value = 1188
$ws.Range("C165").Value = 5
$ws.Range("E165").Value = 15

$ws.Range("D169").Value = 878
$ws.Range("E169").Value = 15

$ws.Range("A204").Value = 'Santa Lucia'

$ws.Range("A205").Value = 'Timor Oriental'

$ws.Range("A214").Value = 'Montserrat'
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = 'Islas Malvinas'
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

